$wb = $excel.ActiveWorkbook

# --- Sheet "initial": add flat-start V/T values (columns C and D) for the
#     PV/PQ buses that did not yet have them, ahead of FDLF Q-limit work ---
$wsInitial = $wb.Worksheets.Item("initial")

$wsInitial.Range("C2").Value = 1
$wsInitial.Range("D2").Value = 0

$wsInitial.Range("C3").Value = 1
$wsInitial.Range("D3").Value = 0

$wsInitial.Range("C5").Value = 1
$wsInitial.Range("D5").Value = 0

$wsInitial.Range("C6").Value = 1
$wsInitial.Range("D6").Value = 0

# Restore selection on "line_imp" before switching away, then make
# "initial" the active/selected sheet with its own cursor position.
$wsLineImp = $wb.Worksheets.Item("line_imp")
$wsLineImp.Activate()
[void]$wsLineImp.Range("D4").Select()

$wsInitial.Activate()
[void]$wsInitial.Range("B3").Select()
